$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F5").WrapText = $true
